$d = $word.ActiveDocument

# --- Part 1: remove the stray _GoBack bookmark around "9 March 2019" ---
$d.Bookmarks("_GoBack").Delete()

# --- Part 2: "VS3" -> "VST3" (keeping the yellow highlight formatting),
#     re-adding a (now collapsed) _GoBack bookmark right after the new "T",
#     matching the run layout Word itself leaves after an in-place edit:
#       run("VS") + run("T") + bookmarkStart + bookmarkEnd + run("3")
$r = $d.Content
$null = $r.Find.Execute("VS3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

# Split "VS" | "3" with a temporary bookmark first; this does not touch the
# run text, so both halves keep the original run's formatting/identity.
$splitPoint = $d.Range($r.Start + 2, $r.Start + 2)
$d.Bookmarks.Add("ZZZTempSplit", $splitPoint)

# Turn the trailing "3" run into "T3" (only this run's text is edited).
$threeRun = $d.Range($r.Start + 2, $r.Start + 3)
$threeRun.Text = "T3"

# Drop the helper bookmark now that "VS" / "T3" are separate runs.
$d.Bookmarks("ZZZTempSplit").Delete()

# Finally split "T3" into "T" | "3" with the real _GoBack bookmark
# (collapsed, right between the two characters).
$goBackPoint = $d.Range($r.Start + 3, $r.Start + 3)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
